$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.865.66"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "1.887.35"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7696"
$ws.Range("E5").Value = "  -1.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.75"
$ws.Range("E6").Value = "  -0.85%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3118"
$ws.Range("E8").Value = "  -0.95%  "
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07176"
$ws.Range("E10").Value = "  -2.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08610"
$ws.Range("E11").Value = "  +6.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7628"
$ws.Range("E12").Value = "  -1.32%  "
$ws.Range("D13").Value = "1.925.20"
$ws.Range("E13").Value = "  +4.58%  "
$ws.Range("E14").Value = "  -2.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.57"
$ws.Range("E15").Value = "  +1.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.161"
$ws.Range("E16").Value = "  -2.60%  "
$ws.Range("D17").Value = "29.909.25"
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.76"
$ws.Range("E18").Value = "  -1.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.45"
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007806"
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("D21").Value = "2.204.08"
$ws.Range("E21").Value = "  +3.53%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.023"
$ws.Range("E23").Value = "  -1.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1647"
$ws.Range("E25").Value = "  +3.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.363"
$ws.Range("E26").Value = "  -1.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.29"
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.74"
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.037"
$ws.Range("E29").Value = "  -0.61%  "
$ws.Range("E30").Value = "  +1.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.532"
$ws.Range("E31").Value = "  -1.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.500"
$ws.Range("E32").Value = "  +0.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.100"
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("E34").Value = "  -1.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.239"
$ws.Range("E35").Value = "  -1.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7447"
$ws.Range("E36").Value = "  -1.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.002"
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.704"
$ws.Range("E38").Value = "  +2.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01957"
$ws.Range("E39").Value = "  +1.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.780"
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4464"
$ws.Range("E41").Value = "  +0.52%  "
$ws.Range("D42").Value = "1.108.36"
$ws.Range("E42").Value = "  -4.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "73.04"
$ws.Range("E43").Value = "  -1.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.070"
$ws.Range("E44").Value = "  +2.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8499"
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.32"
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.658"
$ws.Range("E48").Value = "  +2.10%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.862"
$ws.Range("E49").Value = "  -1.84%  "
$ws.Range("D50").Value = "2.099.92"
$ws.Range("E50").Value = "  +2.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.984"
$ws.Range("E51").Value = "  -3.41%  "
